$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$rng = $hdr.Range
Write-Output ("before collapse: start=" + $rng.Start + " end=" + $rng.End)
$rng.Collapse(1)
Write-Output ("after collapse: start=" + $rng.Start + " end=" + $rng.End)
